$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Nokia"
$ws.Range("C2").Value = "charger"
$ws.Range("D2").Value = "pin head"
$ws.Range("E2").Value = 1000.0
$ws.Range("F2").Value = 50.0
$ws.Range("G2").Value = "27vl"

# Row 3
$ws.Range("C3").Value = "phone"

# Row 4
$ws.Range("C4").Value = "phone"

# Row 5
$ws.Range("C5").Value = "phone"

# Row 6
$ws.Range("C6").Value = "phone"

# Row 7
$ws.Range("B7").Value = "Beat soudio"
$ws.Range("C7").Value = "headset"
$ws.Range("D7").Value = "pin head"
$ws.Range("E7").Value = 50000.0
$ws.Range("G7").Value = "4GB"

# Row 8
$ws.Range("B8").Value = "Kings"
$ws.Range("C8").Value = "headset"
$ws.Range("D8").Value = "pin head"
$ws.Range("E8").Value = 6000.0
$ws.Range("G8").Value = "BT"

# Row 9
$ws.Range("B9").Value = "Kings"
$ws.Range("C9").Value = "headset"
$ws.Range("D9").Value = "USB"
$ws.Range("E9").Value = 12000.0
$ws.Range("G9").Value = "4GB"

# Row 10
$ws.Range("B10").Value = "iphone"
$ws.Range("C10").Value = "charger"
$ws.Range("D10").Value = "2pace"
$ws.Range("E10").Value = 20000.0
$ws.Range("F10").Value = 7.0
$ws.Range("G10").Value = "26VL"

# Row 11
$ws.Range("B11").Value = "tecno"
$ws.Range("C11").Value = "charger"
$ws.Range("D11").Value = "USB"
$ws.Range("E11").Value = 2000.0
$ws.Range("F11").Value = 40.0
$ws.Range("G11").Value = "24VL"

# Row 12
$ws.Range("C12").Value = "phone"
$ws.Range("D12").Value = "black"

# Row 13
$ws.Range("B13").Value = "Beat soudio"
$ws.Range("C13").Value = "headset"
$ws.Range("D13").Value = "USB"
$ws.Range("F13").Value = 7.0
